$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Treatment" query text (B5): CONCAT(REPLACE(...)) -> REPLACE(...)
$newTreatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs003432' AND dgn.anatomic_site = 'C42.1 : Bone marrow'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

$ws.Range("B5").Value = $newTreatmentQuery

# Re-apply the cell's font formatting (12pt Calibri, wrap text). Bumping the
# size away and back mirrors the original author's re-save of this cell,
# which is what gives it its own distinct style entry (separate from the
# other query cells in B6:B7 that keep sharing their original style).
$ws.Range("B5").Font.Size = 13
$ws.Range("B5").Font.Size = 12
$ws.Range("B5").WrapText = $true

# Restore the scroll position / active selection to B5 (top-left cell A4).
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
